$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.047125
$ws.Range("H2").Value = 0.141375
$ws.Range("I2").Value = 0.1108387998127795
$ws.Range("J2").Value = 0.1108387998127795
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3427866666666666
$ws.Range("N2").Value = 1.02836
$ws.Range("O2").Value = 0.9518335355734316
$ws.Range("P2").Value = 0.9518335355734318
$ws.Range("Q2").Value = 0.01615382166666667
$ws.Range("R2").Value = 0.145384395
$ws.Range("S2").Value = 0.1055000867045137
$ws.Range("T2").Value = 0.1055000867045137
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.047125
$ws.Range("H3").Value = 0.141375
$ws.Range("I3").Value = 0.1108387998127795
$ws.Range("J3").Value = 0.1108387998127795
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01734633333333334
$ws.Range("N3").Value = 0.052039
$ws.Range("O3").Value = 0.04816646442656834
$ws.Range("P3").Value = 0.04816646442656834
$ws.Range("Q3").Value = 0.0008174459583333335
$ws.Range("R3").Value = 0.007357013625
$ws.Range("S3").Value = 0.005338713108265771
$ws.Range("T3").Value = 0.005338713108265771
$ws.Range("I4").Value = 0.5286385506557817
$ws.Range("J4").Value = 0.5286385506557816
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3427866666666666
$ws.Range("N4").Value = 1.02836
$ws.Range("O4").Value = 0.9518335355734316
$ws.Range("P4").Value = 0.9518335355734318
$ws.Range("Q4").Value = 0.07704461693777778
$ws.Range("R4").Value = 0.69340155244
$ws.Range("S4").Value = 0.5031759007111073
$ws.Range("T4").Value = 0.5031759007111073
$ws.Range("I5").Value = 0.5286385506557817
$ws.Range("J5").Value = 0.5286385506557816
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01734633333333334
$ws.Range("N5").Value = 0.052039
$ws.Range("O5").Value = 0.04816646442656834
$ws.Range("P5").Value = 0.04816646442656834
$ws.Range("Q5").Value = 0.00389875609788889
$ws.Range("R5").Value = 0.035088804881
$ws.Range("S5").Value = 0.02546264994467436
$ws.Range("T5").Value = 0.02546264994467435
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1532823333333333
$ws.Range("H6").Value = 0.459847
$ws.Range("I6").Value = 0.360522649531439
$ws.Range("J6").Value = 0.360522649531439
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3427866666666666
$ws.Range("N6").Value = 1.02836
$ws.Range("O6").Value = 0.9518335355734316
$ws.Range("P6").Value = 0.9518335355734318
$ws.Range("Q6").Value = 0.05254314010222221
$ws.Range("R6").Value = 0.4728882609199999
$ws.Range("S6").Value = 0.3431575481578108
$ws.Range("T6").Value = 0.3431575481578108
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1532823333333333
$ws.Range("H7").Value = 0.459847
$ws.Range("I7").Value = 0.360522649531439
$ws.Range("J7").Value = 0.360522649531439
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01734633333333334
$ws.Range("N7").Value = 0.052039
$ws.Range("O7").Value = 0.04816646442656834
$ws.Range("P7").Value = 0.04816646442656834
$ws.Range("Q7").Value = 0.002658886448111111
$ws.Range("R7").Value = 0.023929978033
$ws.Range("S7").Value = 0.01736510137362822
$ws.Range("T7").Value = 0.01736510137362822
